$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
# A8: "Volume 32   Number  32" -> "Volume 32   Number  34"
$ws.Range("A8").Value = "Volume 32   Number  34"
# C9: "Report Covering the Week  8/4/2025  Through  8/10/2025"
#     -> "Report Covering the Week  8/18/2025  Through  8/24/2025"
$ws.Range("C9").Value = "Report Covering the Week  8/18/2025  Through  8/24/2025"

# --- Crime Complaints table (rows 15-33): new weekly figures ---
# Row 14 (Murder) is untouched by this update and is used below purely as a
# formatting donor (PasteSpecial -> xlPasteFormats) so that cells which
# switch between a numeric value and the special "0" / "***.*" placeholder
# text keep the workbook's existing cell style (13 = text marker style,
# 14 = integer style) instead of Excel auto-generating a brand new style.

$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 175
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = -26.785714285714
$ws.Range("L16").Value = -28.070175438596
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 133.333333333333
$ws.Range("I17").Value = 169
$ws.Range("J17").Value = 145
$ws.Range("K17").Value = 16.551724137931
$ws.Range("L17").Value = 4.320987654320
$ws.Range("C18").Value = 2
$ws.Range("J14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = -36.538461538461
$ws.Range("L18").Value = -32.653061224489
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -75
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -48
$ws.Range("I19").Value = 139
$ws.Range("J19").Value = 176
$ws.Range("K19").Value = -21.022727272727
$ws.Range("L19").Value = -31.527093596059
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -57.894736842105
$ws.Range("I20").Value = 111
$ws.Range("J20").Value = 116
$ws.Range("K20").Value = -4.310344827586
$ws.Range("L20").Value = -9.016393442622
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -63.636363636363
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = -20.588235294117
$ws.Range("I21").Value = 515
$ws.Range("J21").Value = 555
$ws.Range("K21").Value = -7.207207207207
$ws.Range("L21").Value = -14.309484193011
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 207.142857142857
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 110
$ws.Range("I24").Value = 383
$ws.Range("J24").Value = 376
$ws.Range("K24").Value = 1.861702127659
$ws.Range("L24").Value = 3.234501347708
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = -46.666666666666
$ws.Range("I25").Value = 68
$ws.Range("J25").Value = 106
$ws.Range("K25").Value = -35.849056603773
$ws.Range("L25").Value = -18.072289156626
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 11.111111111111
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 24.242424242424
$ws.Range("I26").Value = 309
$ws.Range("J26").Value = 296
$ws.Range("K26").Value = 4.391891891891
$ws.Range("L26").Value = 28.75
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 80
$ws.Range("D28").Value = 3
$ws.Range("F28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -14.285714285714
$ws.Range("G29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("G30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("F33").Value = 1
$ws.Range("J14").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("I33").Value = 3
